# Weekly fruit/vegetable price update: insert this week's two new price
# rows (Primera / Segunda quality, week of 2021-09-09) above the existing
# historical rows, pushing the old rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 116-117; this shifts the previous rows
# 116/117/118 down to 118/119/120 (with their formatting/values intact).
$ws.Rows("116:117").Insert()

# New row 116: Apio, Americana (o), Primera - week of 2021-09-09
$ws.Cells.Item(116, 1).Value = 9
$ws.Cells.Item(116, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(116, 3).Value = "Metropolitana"
$ws.Cells.Item(116, 4).Value = 44448
$ws.Cells.Item(116, 5).Value = 13
$ws.Cells.Item(116, 6).Value = 100112017
$ws.Cells.Item(116, 7).Value = "Apio"
$ws.Cells.Item(116, 8).Value = "Americana (o)"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 61
$ws.Cells.Item(116, 11).Value = 8000
$ws.Cells.Item(116, 12).Value = 9000
$ws.Cells.Item(116, 13).Value = 8492
$ws.Cells.Item(116, 14).Value = "$/docena de matas"
$ws.Cells.Item(116, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(116, 16).Value = 1415
$ws.Cells.Item(116, 17).Value = 6
$ws.Cells.Item(116, 18).Value = "Hortaliza"

# New row 117: Apio, Americana (o), Segunda - week of 2021-09-09
$ws.Cells.Item(117, 1).Value = 9
$ws.Cells.Item(117, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(117, 3).Value = "Metropolitana"
$ws.Cells.Item(117, 4).Value = 44448
$ws.Cells.Item(117, 5).Value = 13
$ws.Cells.Item(117, 6).Value = 100112017
$ws.Cells.Item(117, 7).Value = "Apio"
$ws.Cells.Item(117, 8).Value = "Americana (o)"
$ws.Cells.Item(117, 9).Value = "Segunda"
$ws.Cells.Item(117, 10).Value = 25
$ws.Cells.Item(117, 11).Value = 6000
$ws.Cells.Item(117, 12).Value = 7000
$ws.Cells.Item(117, 13).Value = 6480
$ws.Cells.Item(117, 14).Value = "$/docena de matas"
$ws.Cells.Item(117, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(117, 16).Value = 1080
$ws.Cells.Item(117, 17).Value = 6
$ws.Cells.Item(117, 18).Value = "Hortaliza"
